$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values would otherwise be
# auto-detected as numbers by Excel (the source data stores prices as text).
$textCells = @("D4", "D5", "D6", "D11", "D12", "D13", "D14", "D17", "D19", "D20", "D21", "D23", "D24", "D25", "D26", "D28", "D31", "D34", "D36", "D37", "D38", "D40", "D41", "D43", "D45", "D46", "D47", "D48", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "64.057.96"
$ws.Range("E2").Value = "  -3.20%  "
$ws.Range("D3").Value = "3.124.91"
$ws.Range("E3").Value = "  -2.59%  "
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").Value = "605.89"
$ws.Range("E5").Value = "  -0.27%  "
$ws.Range("D6").Value = "146.97"
$ws.Range("E6").Value = "  -5.13%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "3.119.67"
$ws.Range("E9").Value = "  -3.72%  "
$ws.Range("E10").Value = "  -5.30%  "
$ws.Range("D11").Value = "5.56"
$ws.Range("E11").Value = "  -2.48%  "
$ws.Range("D12").Value = "0.473"
$ws.Range("E12").Value = "  -5.28%  "
$ws.Range("D13").Value = "0.0000255"
$ws.Range("E13").Value = "  -4.14%  "
$ws.Range("D14").Value = "36.31"
$ws.Range("E14").Value = "  -4.83%  "
$ws.Range("D15").Value = "3.630.78"
$ws.Range("E15").Value = "  -2.70%  "
$ws.Range("D16").Value = "63.946.78"
$ws.Range("E16").Value = "  -3.49%  "
$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D17").Value = "0.113"
$ws.Range("E17").Value = "  +0.06%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.131.21"
$ws.Range("E18").Value = "  -2.49%  "
$ws.Range("D19").Value = "6.91"
$ws.Range("E19").Value = "  -4.45%  "
$ws.Range("D20").Value = "478.18"
$ws.Range("E20").Value = "  -5.35%  "
$ws.Range("D21").Value = "14.47"
$ws.Range("E21").Value = "  -4.60%  "
$ws.Range("E22").Value = "  -3.09%  "
$ws.Range("D23").Value = "7.67"
$ws.Range("E23").Value = "  -3.76%  "
$ws.Range("D24").Value = "13.66"
$ws.Range("E24").Value = "  -5.67%  "
$ws.Range("D25").Value = "83.06"
$ws.Range("E25").Value = "  -2.00%  "
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("E27").Value = "  -2.69%  "
$ws.Range("D28").Value = "8.44"
$ws.Range("E28").Value = "  -5.58%  "
$ws.Range("E29").Value = "  -5.07%  "
$ws.Range("E30").Value = "  -20.88%  "
$ws.Range("D31").Value = "6.82"
$ws.Range("E31").Value = "  -1.27%  "
$ws.Range("E32").Value = "  -0.06%  "
$ws.Range("E33").Value = "  -5.80%  "
$ws.Range("D34").Value = "26.49"
$ws.Range("E34").Value = "  -5.80%  "
$ws.Range("E35").Value = "  -5.45%  "
$ws.Range("D36").Value = "6.04"
$ws.Range("E36").Value = "  -4.96%  "
$ws.Range("D37").Value = "54.42"
$ws.Range("E37").Value = "  -1.62%  "
$ws.Range("D38").Value = "3.11"
$ws.Range("E38").Value = "  +3.91%  "
$ws.Range("D39").Value = "0.0₃0717"
$ws.Range("E39").Value = "  -5.91%  "
$ws.Range("D40").Value = "449.72"
$ws.Range("E40").Value = "  -9.56%  "
$ws.Range("D41").Value = "0.0395"
$ws.Range("E41").Value = "  -5.07%  "
$ws.Range("E42").Value = "  -5.83%  "
$ws.Range("D43").Value = "8.38"
$ws.Range("E43").Value = "  -3.58%  "
$ws.Range("D44").Value = "2.852.19"
$ws.Range("E44").Value = "  -2.25%  "
$ws.Range("D45").Value = "0.269"
$ws.Range("E45").Value = "  -8.22%  "
$ws.Range("D46").Value = "2.25"
$ws.Range("E46").Value = "  -6.95%  "
$ws.Range("D47").Value = "26.43"
$ws.Range("E47").Value = "  -5.21%  "
$ws.Range("D48").Value = "0.998"
$ws.Range("E48").Value = "  -0.05%  "
$ws.Range("E49").Value = "  -2.57%  "
$ws.Range("E50").Value = "  -3.45%  "
$ws.Range("D51").Value = "118.64"
$ws.Range("E51").Value = "  -2.62%  "
